$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected (legacy password hash) which blocks direct
# Range.Value writes via COM. Temporarily lift protection, make the
# required edits, then re-apply sheet protection.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure banner (shared string).
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-04 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) / Percent Change (E) columns for rows 2-8.
$ws.Range("D2").Value = 0.4940369917197517
$ws.Range("E2").Value = 0.001695030171537049

$ws.Range("D3").Value = 0.2482426335251187
$ws.Range("E3").Value = -0.01326030556356295

$ws.Range("D4").Value = 0.09787563012543414
$ws.Range("E4").Value = -0.01053534239862797

$ws.Range("D5").Value = 0.1026300118131569
$ws.Range("E5").Value = 0.001947148817802447

$ws.Range("D6").Value = 0.02974985683184722
$ws.Range("E6").Value = 0.0009596008060648309

$ws.Range("D7").Value = 0.02746487598469143
$ws.Range("E7").Value = -0.01074777040932995

$ws.Range("D8").Value = 0.9999999999999999
$ws.Range("E8").Value = -0.003552321132114122

# Restore sheet protection.
$ws.Protect()
